# Refresh the cryptos price/volume columns (D = Price, E = Volume(1h))
# with the latest scraped figures.
#
# Values that look like plain numbers (e.g. "319.93") get typed with a
# leading apostrophe so Excel stores them as exact text (matching the
# workbook's existing text-only Price/Volume columns) instead of silently
# coercing them to floating point and losing precision/trailing zeros;
# the style is then reset to "Normal" so no stray text-format/quote-prefix
# styling sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '43.278.35'
$ws.Cells.Item(2,5).Value = '  -4.60%  '

$ws.Cells.Item(3,4).Value = '2.239.56'
$ws.Cells.Item(3,5).Value = '  -5.52%  '

$ws.Cells.Item(4,5).Value = '  -0.10%  '

$ws.Cells.Item(5,4).Value = '''319.93'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  +1.76%  '

$ws.Cells.Item(6,4).Value = '''101.33'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  -6.29%  '

$ws.Cells.Item(7,4).Value = '''0.588'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  -7.38%  '

$ws.Cells.Item(8,5).Value = '  -0.14%  '

$ws.Cells.Item(9,4).Value = '''0.565'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  -7.79%  '

$ws.Cells.Item(10,4).Value = '''37.16'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  -8.97%  '

$ws.Cells.Item(11,4).Value = '''54.49'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  -2.56%  '

$ws.Cells.Item(12,4).Value = '''0.0827'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  -9.73%  '

$ws.Cells.Item(13,4).Value = '''7.72'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  -9.02%  '

$ws.Cells.Item(14,5).Value = '  -1.19%  '

$ws.Cells.Item(15,4).Value = '2.575.15'

$ws.Cells.Item(16,4).Value = '''0.865'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  -11.82%  '

$ws.Cells.Item(17,4).Value = '''14.42'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = '  -6.02%  '

$ws.Cells.Item(18,4).Value = '2.236.00'
$ws.Cells.Item(18,5).Value = '  -5.48%  '

$ws.Cells.Item(19,4).Value = '43.157.21'
$ws.Cells.Item(19,5).Value = '  -4.74%  '

$ws.Cells.Item(20,4).Value = '''14.46'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  -6.16%  '

$ws.Cells.Item(21,4).Value = '0.0₃0969'
$ws.Cells.Item(21,5).Value = '  -8.59%  '

$ws.Cells.Item(22,4).Value = '''6.56'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  -9.90%  '

$ws.Cells.Item(23,4).Value = '''65.57'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  -10.50%  '

$ws.Cells.Item(24,4).Value = '''3.19'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  -11.53%  '

$ws.Cells.Item(25,4).Value = '''238.78'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  -8.55%  '

$ws.Cells.Item(26,4).Value = '''2.17'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -7.96%  '

$ws.Cells.Item(27,5).Value = '  -0.45%  '

$ws.Cells.Item(28,5).Value = '  +2.17%  '

$ws.Cells.Item(29,4).Value = '''10.09'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  -9.66%  '

$ws.Cells.Item(30,5).Value = '  -2.34%  '

$ws.Cells.Item(31,4).Value = '''6.43'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  -14.10%  '

$ws.Cells.Item(32,4).Value = '''35.49'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -4.21%  '

$ws.Cells.Item(33,4).Value = '''20.53'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  -7.97%  '

$ws.Cells.Item(34,4).Value = '''0.0880'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  -8.91%  '

$ws.Cells.Item(35,4).Value = '''153.67'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -7.74%  '

$ws.Cells.Item(36,5).Value = '  -4.24%  '

$ws.Cells.Item(37,4).Value = '''3.18'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  +9.14%  '

$ws.Cells.Item(38,4).Value = '''1.99'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  +5.61%  '

$ws.Cells.Item(39,5).Value = '  -6.83%  '

$ws.Cells.Item(40,4).Value = '''4.48'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  -4.45%  '

$ws.Cells.Item(41,5).Value = '  -10.79%  '

$ws.Cells.Item(42,4).Value = '''3.71'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  -6.43%  '

$ws.Cells.Item(43,4).Value = '''0.0323'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = '  -8.55%  '

$ws.Cells.Item(44,4).Value = '''12.94'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = '  -1.01%  '

$ws.Cells.Item(45,5).Value = '  -0.06%  '

$ws.Cells.Item(46,4).Value = '1.798.98'
$ws.Cells.Item(46,5).Value = '  -1.15%  '

$ws.Cells.Item(47,4).Value = '''87.41'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  -11.31%  '

$ws.Cells.Item(48,5).Value = '  -9.44%  '

$ws.Cells.Item(49,4).Value = '''76.95'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  -7.32%  '

$ws.Cells.Item(50,4).Value = '''5.34'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = '  -10.02%  '

$ws.Cells.Item(51,4).Value = '''59.61'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -14.93%  '

